$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) "Ativacao:" date text 01/01/2020 -> 01/01/2023 -------------------
# The cells store the date as literal text (shared string), not a real
# Excel date. Writing the string straight into .Value would make Excel's
# COM layer auto-convert it into a date serial (and pull in a new number
# format/style), so instead we park the literal text in a scratch cell via
# a formula (whose cached result is plain text), copy it, and paste-special
# *values only* into the target cells. That preserves both the text type
# and the existing per-cell style.
$ws.Range("E1").Formula = '="01/01/2023"'
$ws.Range("E1").Copy()
$ws.Range("B8").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("C8").PasteSpecial(-4163)
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("C15").PasteSpecial(-4163)
$ws.Range("E1").Value = ""

# --- 2) Fill in the new English descriptions ------------------------------
# Columns B/C were previously blank on these rows. Column C's width/style
# definition is unambiguous, so a plain .Value assignment already lands on
# the right style (s=3). Column B has two overlapping <col> definitions, so
# a brand-new cell there would otherwise inherit the wrong style; copy the
# number/alignment formatting from an already-populated column-B cell first.
$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("B16").PasteSpecial(-4122)

$ws.Range("B11").Value = "To present concepts about renewable sources for the generation of thermal, electrical and vehicular energy, among others"
$ws.Range("C11").Value = "To present concepts about renewable sources for the generation of thermal, electrical and vehicular energy, among others"

$ws.Range("B14").Value = "Renewable sources and clean technologies for energy generation. Study of current national and global systems."
$ws.Range("C14").Value = "Renewable sources and clean technologies for energy generation. Study of current national and global systems."

$ws.Range("B16").Value = "National and global energy systems: renewable and fossil sources. Energy generation from renewable sources: solar thermal and photovoltaic; wind; maritime. Generation of biomass for energy purposes. Management of urban solid waste: recyclable and non-recyclable; enterprise programs for reverse logistics; the issue of polymers; reforestation; processing of domestic wet waste. Integration of renewable sources for energy generation: hybrid thermal cycles"
$ws.Range("C16").Value = "National and global energy systems: renewable and fossil sources. Energy generation from renewable sources: solar thermal and photovoltaic; wind; maritime. Generation of biomass for energy purposes. Management of urban solid waste: recyclable and non-recyclable; enterprise programs for reverse logistics; the issue of polymers; reforestation; processing of domestic wet waste. Integration of renewable sources for energy generation: hybrid thermal cycles"
